$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell H1 "Letra", matching the centered style used by the other headers (e.g. G1)
$ws.Range("H1").Value = "Letra"
$ws.Range("H1").HorizontalAlignment = -4108
$ws.Range("H1").VerticalAlignment = -4108

# Update E2 value from 73 to 16
$ws.Range("E2").Value = 16

# Add new cell H2 "Q"
$ws.Range("H2").Value = "Q"

# Update the active selection to H6 (as recorded in the saved view state)
$ws.Range("H6").Select()
